$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Materialize rows 1105:1126 first (copying the blank styled template from row 1090, which is
#    about to be overwritten with data) so every new row - both the ones that get data (1105:1108)
#    and the fresh blank trailer rows (1109:1126) - starts from the correct "blank" cell styles
#    (s="3" text cells / s="4" formula cells) instead of whatever the bare column defaults are.
$ws.Range("A1090:J1090").Copy($ws.Range("A1105:J1126"))

# 2) Fill in the new timesheet entries for 2024-06-16 .. 2024-06-24 (rows 1090:1108).
$ws.Range("A1090").Value = "2024-06-16"
$ws.Range("B1090").Value = "13:45"
$ws.Range("C1090").Value = "14:45"
$ws.Range("D1090").Value = "1h 00m"
$ws.Range("E1090").Value = "#maintenance"
$ws.Range("G1090").Value = "'False"
$ws.Range("H1090").Value = "'False"

$ws.Range("A1091").Value = "2024-06-16"
$ws.Range("B1091").Value = "14:45"
$ws.Range("C1091").Value = "20:15"
$ws.Range("D1091").Value = "5h 30m"
$ws.Range("E1091").Value = "#python"
$ws.Range("F1091").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1091").Value = "'True"
$ws.Range("H1091").Value = "'False"

$ws.Range("A1092").Value = "2024-06-17"
$ws.Range("B1092").Value = "10:30"
$ws.Range("C1092").Value = "11:00"
$ws.Range("D1092").Value = "0h 30m"
$ws.Range("E1092").Value = "#python"
$ws.Range("F1092").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1092").Value = "'True"
$ws.Range("H1092").Value = "'False"

$ws.Range("A1093").Value = "2024-06-17"
$ws.Range("B1093").Value = "12:00"
$ws.Range("C1093").Value = "14:00"
$ws.Range("D1093").Value = "2h 00m"
$ws.Range("E1093").Value = "#python"
$ws.Range("F1093").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1093").Value = "'True"
$ws.Range("H1093").Value = "'False"

$ws.Range("A1094").Value = "2024-06-17"
$ws.Range("B1094").Value = "16:30"
$ws.Range("C1094").Value = "18:00"
$ws.Range("D1094").Value = "1h 30m"
$ws.Range("E1094").Value = "#python"
$ws.Range("F1094").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1094").Value = "'True"
$ws.Range("H1094").Value = "'False"

$ws.Range("A1095").Value = "2024-06-17"
$ws.Range("B1095").Value = "20:00"
$ws.Range("C1095").Value = "21:15"
$ws.Range("D1095").Value = "1h 15m"
$ws.Range("E1095").Value = "#python"
$ws.Range("F1095").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1095").Value = "'True"
$ws.Range("H1095").Value = "'False"

$ws.Range("A1096").Value = "2024-06-17"
$ws.Range("B1096").Value = "22:00"
$ws.Range("C1096").Value = "22:15"
$ws.Range("D1096").Value = "0h 15m"
$ws.Range("E1096").Value = "#python"
$ws.Range("F1096").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1096").Value = "'True"
$ws.Range("H1096").Value = "'False"

$ws.Range("A1097").Value = "2024-06-18"
$ws.Range("B1097").Value = "10:00"
$ws.Range("C1097").Value = "14:00"
$ws.Range("D1097").Value = "4h 00m"
$ws.Range("E1097").Value = "#maintenance"
$ws.Range("G1097").Value = "'False"
$ws.Range("H1097").Value = "'False"

$ws.Range("A1098").Value = "2024-06-18"
$ws.Range("B1098").Value = "15:15"
$ws.Range("C1098").Value = "17:30"
$ws.Range("D1098").Value = "2h 15m"
$ws.Range("E1098").Value = "#python"
$ws.Range("F1098").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1098").Value = "'True"
$ws.Range("H1098").Value = "'False"

$ws.Range("A1099").Value = "2024-06-18"
$ws.Range("B1099").Value = "19:30"
$ws.Range("C1099").Value = "22:00"
$ws.Range("D1099").Value = "2h 30m"
$ws.Range("E1099").Value = "#python"
$ws.Range("F1099").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1099").Value = "'True"
$ws.Range("H1099").Value = "'False"

$ws.Range("A1100").Value = "2024-06-20"
$ws.Range("B1100").Value = "08:15"
$ws.Range("C1100").Value = "08:45"
$ws.Range("D1100").Value = "0h 30m"
$ws.Range("E1100").Value = "#python"
$ws.Range("F1100").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1100").Value = "'True"
$ws.Range("H1100").Value = "'False"

$ws.Range("A1101").Value = "2024-06-20"
$ws.Range("B1101").Value = "17:00"
$ws.Range("C1101").Value = "17:45"
$ws.Range("D1101").Value = "0h 45m"
$ws.Range("E1101").Value = "#python"
$ws.Range("F1101").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1101").Value = "'True"
$ws.Range("H1101").Value = "'False"

$ws.Range("A1102").Value = "2024-06-21"
$ws.Range("B1102").Value = "08:00"
$ws.Range("C1102").Value = "08:45"
$ws.Range("D1102").Value = "0h 45m"
$ws.Range("E1102").Value = "#python"
$ws.Range("F1102").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1102").Value = "'True"
$ws.Range("H1102").Value = "'False"

$ws.Range("A1103").Value = "2024-06-23"
$ws.Range("B1103").Value = "17:30"
$ws.Range("C1103").Value = "19:30"
$ws.Range("D1103").Value = "2h 00m"
$ws.Range("E1103").Value = "#python"
$ws.Range("F1103").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1103").Value = "'True"
$ws.Range("H1103").Value = "'False"

$ws.Range("A1104").Value = "2024-06-23"
$ws.Range("B1104").Value = "20:30"
$ws.Range("C1104").Value = "22:15"
$ws.Range("D1104").Value = "1h 45m"
$ws.Range("E1104").Value = "#python"
$ws.Range("F1104").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1104").Value = "'True"
$ws.Range("H1104").Value = "'False"

$ws.Range("A1105").Value = "2024-06-23"
$ws.Range("B1105").Value = "22:30"
$ws.Range("C1105").Value = "00:00"
$ws.Range("D1105").Value = "1h 30m"
$ws.Range("E1105").Value = "#python"
$ws.Range("F1105").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1105").Value = "'True"
$ws.Range("H1105").Value = "'False"

$ws.Range("A1106").Value = "2024-06-24"
$ws.Range("B1106").Value = "08:30"
$ws.Range("C1106").Value = "10:30"
$ws.Range("D1106").Value = "2h 00m"
$ws.Range("E1106").Value = "#python"
$ws.Range("F1106").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1106").Value = "'True"
$ws.Range("H1106").Value = "'True"

$ws.Range("A1107").Value = "2024-06-24"
$ws.Range("B1107").Value = "11:45"
$ws.Range("C1107").Value = "12:15"
$ws.Range("D1107").Value = "0h 30m"
$ws.Range("E1107").Value = "#python"
$ws.Range("F1107").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1107").Value = "'True"
$ws.Range("H1107").Value = "'True"

$ws.Range("A1108").Value = "2024-06-24"
$ws.Range("B1108").Value = "13:45"
$ws.Range("C1108").Value = "14:15"
$ws.Range("D1108").Value = "0h 30m"
$ws.Range("E1108").Value = "#python"
$ws.Range("F1108").Value = "nwtraderaanalytics v4.0.0"
$ws.Range("G1108").Value = "'True"
$ws.Range("H1108").Value = "'True"

# 3) YEAR()/MONTH() helper formulas for the new rows, written per contiguous block so the engine
#    groups them the same way Excel would when the formula is filled down a block at a time.
$ws.Range("I1090").Formula = "=YEAR(A1090)"
$ws.Range("J1090").Formula = "=MONTH(A1090)"
$ws.Range("I1091:I1097").Formula = "=YEAR(A1091)"
$ws.Range("J1091:J1097").Formula = "=MONTH(A1091)"
$ws.Range("I1098:I1099").Formula = "=YEAR(A1098)"
$ws.Range("J1098:J1099").Formula = "=MONTH(A1098)"
$ws.Range("I1100:I1102").Formula = "=YEAR(A1100)"
$ws.Range("J1100:J1102").Formula = "=MONTH(A1100)"
$ws.Range("I1103:I1105").Formula = "=YEAR(A1103)"
$ws.Range("J1103:J1105").Formula = "=MONTH(A1103)"
$ws.Range("I1106:I1108").Formula = "=YEAR(A1106)"
$ws.Range("J1106:J1108").Formula = "=MONTH(A1106)"

# 4) Update the view state: scrolled-in frozen pane + the last-active selected cell.
$ws.Range("F1101").Select()

